# Switched to EPPlus instead of Excel.Interop
$wb = $excel.ActiveWorkbook

# Rename the "Process" sheet to "burp"
$wsProcess = $wb.Worksheets.Item("Process")
$wsProcess.Name = "burp"

# Input sheet: D4 changes from 8 to 9, selection moves to E4
$wsInput = $wb.Worksheets.Item("Input")
$wsInput.Range("D4").Value = 9
$wsInput.Range("E4").Select()

# burp (formerly Process) sheet: selection moves to B5
$wsProcess.Range("B5").Select()

# Output sheet: update VLOOKUP formula to reference the renamed sheet
# and add an explicit FALSE for exact match, regroup the parentheses
$wsOutput = $wb.Worksheets.Item("Output")
$wsOutput.Range("C5").Formula = "=((VLOOKUP(Input!D5,burp!A2:B5,2,FALSE))*C3)/1000"

# Keep Output as the active sheet/tab (matches original workbook state)
$wsOutput.Select()
$wsOutput.Range("C5").Select()
